# Scheduled cryptos-list refresh (GitHub Actions cron).
# Coinranking export keeps Price/Volume columns as plain text, even when a
# refreshed value happens to look like a number (e.g. "7.03"). Excel's COM
# layer auto-coerces such assignments to numeric cells, so for those specific
# addresses we briefly force Text number format before writing the value, then
# restore the Normal cell style so formatting matches the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textAddrs = @('D5','D6','D8','D10','D11','D12','D14','D18','D19','D20','D21','D22','D23','D27','D28','D29','D30','D33','D37','D39','D40','D42','D43','D44','D45','D46','D47','D48','D49','D50')
foreach ($addr in $textAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.419.79'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '3.669.71'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '646.69'
$ws.Range("E5").Value = '  -5.00%  '
$ws.Range("D6").Value = '158.91'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.495'
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '7.03'
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '0.0000229'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '4.291.66'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").Value = '32.35'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").Value = '3.668.80'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '69.410.45'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '15.89'
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '464.74'
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").Value = '9.65'
$ws.Range("E21").Value = '  -3.81%  '
$ws.Range("D22").Value = '0.638'
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").Value = '79.34'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("D24").Value = '3.816.22'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").Value = '10.74'
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("D28").Value = '8.83'
$ws.Range("E28").Value = '  -3.42%  '
$ws.Range("D29").Value = '2.60'
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("D30").Value = '1.65'
$ws.Range("E30").Value = '  -6.68%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").Value = '26.53'
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("D35").Value = '3.660.23'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  +2.67%  '
$ws.Range("D37").Value = '8.32'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '5.88'
$ws.Range("E39").Value = '  -6.05%  '
$ws.Range("D40").Value = '179.12'
$ws.Range("E40").Value = '  +4.70%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '2.18'
$ws.Range("E42").Value = '  -3.73%  '
$ws.Range("D43").Value = '0.0889'
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("D44").Value = '0.927'
$ws.Range("E44").Value = '  -1.74%  '
$ws.Range("D45").Value = '46.79'
$ws.Range("E45").Value = '  -1.92%  '
$ws.Range("D46").Value = '2.69'
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '26.89'
$ws.Range("E47").Value = '  -4.39%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '1.25'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = '0.000264'
$ws.Range("E49").Value = '  -5.00%  '
$ws.Range("D50").Value = '7.76'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  -4.54%  '

# Restore default (General) cell style on the cells we had to force to Text,
# so the saved workbook's styling matches an ordinary text-literal cell.
foreach ($addr in $textAddrs) {
    $ws.Range($addr).Style = "Normal"
}
